# Auto-generated edit script: updates market/profit calculation cells
# across multiple worksheets per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: columns H,I,J,K,L,M,N
$ws.Cells.Item(6, 8).Value = 160296.16
$ws.Cells.Item(6, 9).Value = 202808.47
$ws.Cells.Item(6, 10).Value = 875
$ws.Cells.Item(6, 11).Value = 608425.41
$ws.Cells.Item(6, 12).Value = 2625
$ws.Cells.Item(6, 13).Value = -608313.41
$ws.Cells.Item(6, 14).Value = -2849
# Row 76: columns H,I,J,K,L,M,N
$ws.Cells.Item(76, 8).Value = 4549.8335
$ws.Cells.Item(76, 9).Value = 4000
$ws.Cells.Item(76, 10).Value = 4659.8
$ws.Cells.Item(76, 11).Value = 4000
$ws.Cells.Item(76, 12).Value = 4659.8
$ws.Cells.Item(76, 13).Value = -3685
$ws.Cells.Item(76, 14).Value = -5289.8
# Row 79: columns H,I,J,K,L,M,N
$ws.Cells.Item(79, 8).Value = 4549.8335
$ws.Cells.Item(79, 9).Value = 4000
$ws.Cells.Item(79, 10).Value = 4659.8
$ws.Cells.Item(79, 11).Value = 4000
$ws.Cells.Item(79, 12).Value = 4659.8
$ws.Cells.Item(79, 13).Value = -2908
$ws.Cells.Item(79, 14).Value = -6843.8
# Row 113: columns H,I,J,K,L,M,N
$ws.Cells.Item(113, 8).Value = 101860
$ws.Cells.Item(113, 9).Value = 168435
$ws.Cells.Item(113, 10).Value = 1997.5
$ws.Cells.Item(113, 11).Value = 168435
$ws.Cells.Item(113, 12).Value = 1997.5
$ws.Cells.Item(113, 13).Value = -165181
$ws.Cells.Item(113, 14).Value = -8505.5
# Row 137: columns H,I,J,K,L,M,N
$ws.Cells.Item(137, 8).Value = 1507.8334
$ws.Cells.Item(137, 9).Value = 1212.4073
$ws.Cells.Item(137, 10).Value = 4166.6665
$ws.Cells.Item(137, 11).Value = 3637.2219
$ws.Cells.Item(137, 12).Value = 12499.9995
$ws.Cells.Item(137, 13).Value = -1087.2219
$ws.Cells.Item(137, 14).Value = -17599.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 32: columns H,I,J,K,L,M,N
$ws.Cells.Item(32, 8).Value = 32686.16
$ws.Cells.Item(32, 9).Value = 6069.4316
$ws.Cells.Item(32, 10).Value = 130280.836
$ws.Cells.Item(32, 11).Value = 6069.4316
$ws.Cells.Item(32, 12).Value = 130280.836
$ws.Cells.Item(32, 13).Value = -5782.4316
$ws.Cells.Item(32, 14).Value = -130854.836
# Row 64: columns H,J,L,N
$ws.Cells.Item(64, 8).Value = 39993.332
$ws.Cells.Item(64, 10).Value = 39993.332
$ws.Cells.Item(64, 12).Value = 39993.332
$ws.Cells.Item(64, 14).Value = -40489.332
# Row 67: columns H,J,L,N
$ws.Cells.Item(67, 8).Value = 39993.332
$ws.Cells.Item(67, 10).Value = 39993.332
$ws.Cells.Item(67, 12).Value = 39993.332
$ws.Cells.Item(67, 14).Value = -41709.332
# Row 102: columns H,I,J,K,L,M,N
$ws.Cells.Item(102, 8).Value = 57923.668
$ws.Cells.Item(102, 9).Value = 92533.55
$ws.Cells.Item(102, 10).Value = 3536.7144
$ws.Cells.Item(102, 11).Value = 92533.55
$ws.Cells.Item(102, 12).Value = 3536.7144
$ws.Cells.Item(102, 13).Value = -90911.55
$ws.Cells.Item(102, 14).Value = -6780.7144
# Row 122: columns H,I,K,M
$ws.Cells.Item(122, 8).Value = 1204.6923
$ws.Cells.Item(122, 9).Value = 1143.4166
$ws.Cells.Item(122, 11).Value = 3430.2498
$ws.Cells.Item(122, 13).Value = -980.2498000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 62: columns H,J,L,N
$ws.Cells.Item(62, 8).Value = 45794
$ws.Cells.Item(62, 10).Value = 45794
$ws.Cells.Item(62, 12).Value = 45794
$ws.Cells.Item(62, 14).Value = -47166
# Row 65: columns H,J,L,N
$ws.Cells.Item(65, 8).Value = 45794
$ws.Cells.Item(65, 10).Value = 45794
$ws.Cells.Item(65, 12).Value = 137382
$ws.Cells.Item(65, 14).Value = -144246
# Row 105: columns H,I,J,K,L,M,N
$ws.Cells.Item(105, 8).Value = 78841.234
$ws.Cells.Item(105, 9).Value = 78522.30499999999
$ws.Cells.Item(105, 10).Value = 79160.16
$ws.Cells.Item(105, 11).Value = 78522.30499999999
$ws.Cells.Item(105, 12).Value = 79160.16
$ws.Cells.Item(105, 13).Value = -76775.30499999999
$ws.Cells.Item(105, 14).Value = -82654.16
# Row 134: columns H,I,J,K,L,M,N
$ws.Cells.Item(134, 8).Value = 3279.2334
$ws.Cells.Item(134, 9).Value = 3468.44
$ws.Cells.Item(134, 10).Value = 2333.2
$ws.Cells.Item(134, 11).Value = 10405.32
$ws.Cells.Item(134, 12).Value = 6999.599999999999
$ws.Cells.Item(134, 13).Value = -7870.32
$ws.Cells.Item(134, 14).Value = -12069.6

$ws = $wb.Worksheets.Item("CRP")
# Row 22: columns H,I,J,K,L,M,N
$ws.Cells.Item(22, 8).Value = 595.5
$ws.Cells.Item(22, 9).Value = 394.57144
$ws.Cells.Item(22, 10).Value = 2002
$ws.Cells.Item(22, 11).Value = 394.57144
$ws.Cells.Item(22, 12).Value = 2002
$ws.Cells.Item(22, 13).Value = -44.57144
$ws.Cells.Item(22, 14).Value = -2702
# Row 31: columns H,I,J,K,L,M,N
$ws.Cells.Item(31, 8).Value = 40297.12
$ws.Cells.Item(31, 9).Value = 1294.6666
$ws.Cells.Item(31, 10).Value = 81249.7
$ws.Cells.Item(31, 11).Value = 1294.6666
$ws.Cells.Item(31, 12).Value = 81249.7
$ws.Cells.Item(31, 13).Value = -999.6666
$ws.Cells.Item(31, 14).Value = -81839.7
# Row 34: columns H,I,J,K,L,M,N
$ws.Cells.Item(34, 8).Value = 40297.12
$ws.Cells.Item(34, 9).Value = 1294.6666
$ws.Cells.Item(34, 10).Value = 81249.7
$ws.Cells.Item(34, 11).Value = 1294.6666
$ws.Cells.Item(34, 12).Value = 81249.7
$ws.Cells.Item(34, 13).Value = -1092.6666
$ws.Cells.Item(34, 14).Value = -81653.7
# Row 99: columns H,I,K,M
$ws.Cells.Item(99, 8).Value = 18999.143
$ws.Cells.Item(99, 9).Value = 4240
$ws.Cells.Item(99, 11).Value = 4240
$ws.Cells.Item(99, 13).Value = -2742
# Row 126: columns H,I,K,M
$ws.Cells.Item(126, 8).Value = 18999.143
$ws.Cells.Item(126, 9).Value = 4240
$ws.Cells.Item(126, 11).Value = 12720
$ws.Cells.Item(126, 13).Value = -10250
# Row 132: columns H,I,J,K,L,M,N
$ws.Cells.Item(132, 8).Value = 3217.6304
$ws.Cells.Item(132, 9).Value = 3229.0571
$ws.Cells.Item(132, 10).Value = 3181.2727
$ws.Cells.Item(132, 11).Value = 9687.1713
$ws.Cells.Item(132, 12).Value = 9543.8181
$ws.Cells.Item(132, 13).Value = -7157.1713
$ws.Cells.Item(132, 14).Value = -14603.8181
# Row 134: columns H,I,K,M
$ws.Cells.Item(134, 8).Value = 1243
$ws.Cells.Item(134, 9).Value = 1243
$ws.Cells.Item(134, 11).Value = 3729
$ws.Cells.Item(134, 13).Value = -1194

$ws = $wb.Worksheets.Item("CUL")
# Row 33: columns H,J,L,N
$ws.Cells.Item(33, 8).Value = 1022.2381
$ws.Cells.Item(33, 10).Value = 1619.3077
$ws.Cells.Item(33, 12).Value = 9715.8462
$ws.Cells.Item(33, 14).Value = -10281.8462
# Row 131: columns H,I,J,K,L,M,N
$ws.Cells.Item(131, 8).Value = 809.26
$ws.Cells.Item(131, 9).Value = 489.18182
$ws.Cells.Item(131, 10).Value = 848.82025
$ws.Cells.Item(131, 11).Value = 1467.54546
$ws.Cells.Item(131, 12).Value = 2546.46075
$ws.Cells.Item(131, 13).Value = 3572.45454
$ws.Cells.Item(131, 14).Value = -12626.46075

$ws = $wb.Worksheets.Item("GSM")
# Row 80: columns H,J,L,N
$ws.Cells.Item(80, 8).Value = 250251740
$ws.Cells.Item(80, 10).Value = 2000
$ws.Cells.Item(80, 12).Value = 2000
$ws.Cells.Item(80, 14).Value = -3996
# Row 83: columns H,J,L,N
$ws.Cells.Item(83, 8).Value = 250251740
$ws.Cells.Item(83, 10).Value = 2000
$ws.Cells.Item(83, 12).Value = 10000
$ws.Cells.Item(83, 14).Value = -19984
# Row 102: columns H,I,K,M
$ws.Cells.Item(102, 8).Value = 2041.8182
$ws.Cells.Item(102, 9).Value = 2041.8182
$ws.Cells.Item(102, 11).Value = 2041.8182
$ws.Cells.Item(102, 13).Value = -419.8181999999999
# Row 104: columns H,J,L,N
$ws.Cells.Item(104, 8).Value = 43998.332
$ws.Cells.Item(104, 10).Value = 43998.332
$ws.Cells.Item(104, 12).Value = 43998.332
$ws.Cells.Item(104, 14).Value = -50986.332
# Row 122: columns H,I,J,K,L,M,N
$ws.Cells.Item(122, 8).Value = 2254.3076
$ws.Cells.Item(122, 9).Value = 2275.5
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 6826.5
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -4376.5
$ws.Cells.Item(122, 14).Value = -10900
# Row 123: columns H,J,L,N
$ws.Cells.Item(123, 8).Value = 9326
$ws.Cells.Item(123, 10).Value = 9326
$ws.Cells.Item(123, 12).Value = 9326
$ws.Cells.Item(123, 14).Value = -14226
# Row 126: columns H,I,K,M
$ws.Cells.Item(126, 8).Value = 5066
$ws.Cells.Item(126, 9).Value = 5103.7144
$ws.Cells.Item(126, 11).Value = 15311.1432
$ws.Cells.Item(126, 13).Value = -12841.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 7: columns H,I,J,K,L,M,N
$ws.Cells.Item(7, 8).Value = 2281.9092
$ws.Cells.Item(7, 9).Value = 1515.2307
$ws.Cells.Item(7, 10).Value = 3389.3333
$ws.Cells.Item(7, 11).Value = 1515.2307
$ws.Cells.Item(7, 12).Value = 3389.3333
$ws.Cells.Item(7, 13).Value = -1403.2307
$ws.Cells.Item(7, 14).Value = -3613.3333
# Row 40: columns H,I,J,K,L,M,N
$ws.Cells.Item(40, 8).Value = 80849.84
$ws.Cells.Item(40, 9).Value = 206111.6
$ws.Cells.Item(40, 10).Value = 2561.25
$ws.Cells.Item(40, 11).Value = 206111.6
$ws.Cells.Item(40, 12).Value = 2561.25
$ws.Cells.Item(40, 13).Value = -205975.6
$ws.Cells.Item(40, 14).Value = -2833.25
# Row 82: columns H,I,J,K,L,M,N
$ws.Cells.Item(82, 8).Value = 1328.9375
$ws.Cells.Item(82, 9).Value = 845.8
$ws.Cells.Item(82, 10).Value = 1548.5454
$ws.Cells.Item(82, 11).Value = 845.8
$ws.Cells.Item(82, 12).Value = 1548.5454
$ws.Cells.Item(82, 13).Value = -484.8
$ws.Cells.Item(82, 14).Value = -2270.5454
# Row 85: columns H,I,J,K,L,M,N
$ws.Cells.Item(85, 8).Value = 1328.9375
$ws.Cells.Item(85, 9).Value = 845.8
$ws.Cells.Item(85, 10).Value = 1548.5454
$ws.Cells.Item(85, 11).Value = 845.8
$ws.Cells.Item(85, 12).Value = 1548.5454
$ws.Cells.Item(85, 13).Value = 402.2
$ws.Cells.Item(85, 14).Value = -4044.5454
# Row 122: columns H,I,K,M
$ws.Cells.Item(122, 8).Value = 1652
$ws.Cells.Item(122, 9).Value = 1468
$ws.Cells.Item(122, 11).Value = 4404
$ws.Cells.Item(122, 13).Value = -1954
# Row 126: columns H,I,J,K,L,M,N
$ws.Cells.Item(126, 8).Value = 2281.9092
$ws.Cells.Item(126, 9).Value = 1515.2307
$ws.Cells.Item(126, 10).Value = 3389.3333
$ws.Cells.Item(126, 11).Value = 4545.6921
$ws.Cells.Item(126, 12).Value = 10167.9999
$ws.Cells.Item(126, 13).Value = -2075.6921
$ws.Cells.Item(126, 14).Value = -15107.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 118: columns H,J,L,N
$ws.Cells.Item(118, 8).Value = 39980
$ws.Cells.Item(118, 10).Value = 39980
$ws.Cells.Item(118, 12).Value = 39980
$ws.Cells.Item(118, 14).Value = -43294
# Row 126: columns H,I,J,K,L,M,N
$ws.Cells.Item(126, 8).Value = 2086.4
$ws.Cells.Item(126, 9).Value = 1909.1428
$ws.Cells.Item(126, 10).Value = 2500
$ws.Cells.Item(126, 11).Value = 5727.428400000001
$ws.Cells.Item(126, 12).Value = 7500
$ws.Cells.Item(126, 13).Value = -3257.428400000001
$ws.Cells.Item(126, 14).Value = -12440
